# Revert the "User" creation so it no longer allows an empty-string url:
# the row that was added with handle="" / url="https://twitter.com/" is
# removed again (it used to be row 12; account_8 used to be row 13, now
# it moves up to row 12).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 is the spurious "empty handle, bare https://twitter.com/ url"
# entry - delete it outright; this shifts row 13 (account_8) up to 12.
$ws.Rows.Item(12).Delete()

# The engine's Hyperlinks.Delete() only works against the whole-sheet
# collection (per-item Delete() is a no-op here), and row deletion does
# not renumber/drop the hyperlink anchors on its own, so clear every
# hyperlink and rebuild the (now 11, was 12) set against the shifted data.
$ws.Range("A1").Hyperlinks.Delete()

for ($r = 2; $r -le 12; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $url = $cell.Value2
    $ws.Hyperlinks.Add($cell, $url) | Out-Null
    # Hyperlinks.Add() re-stamps cell style; put it back to the shared
    # "Hyperlink" look the column already used so styles.xml stays stable.
    $cell.Style = "Hyperlink"
}

# Move the live selection the way the saved file shows (below the data).
$ws.Range("B15").Select() | Out-Null
